# Auto-generated edit script applying the diff to Golem_Profits workbook (sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2103.5
$ws.Range("I12").Value = 2138
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 2138
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -1968
$ws.Range("N12").Value = -2340
$ws.Range("H34").Value = 3099.5715
$ws.Range("I34").Value = 2782.8333
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2782.8333
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2579.8333
$ws.Range("N34").Value = -5406
$ws.Range("H36").Value = 3099.5715
$ws.Range("I36").Value = 2782.8333
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 2782.8333
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -2067.8333
$ws.Range("N36").Value = -6430
$ws.Range("H47").Value = 25000
$ws.Range("I47").Value = 25000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 25000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -24028
$ws.Range("N47").ClearContents()
$ws.Range("H86").Value = 4133.778
$ws.Range("I86").Value = 3600.5715
$ws.Range("K86").Value = 3600.5715
$ws.Range("M86").Value = -2477.5715
$ws.Range("H89").Value = 4133.778
$ws.Range("I89").Value = 3600.5715
$ws.Range("K89").Value = 18002.8575
$ws.Range("M89").Value = -12386.8575
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3377.6667
$ws.Range("I35").Value = 3253.4
$ws.Range("K35").Value = 3253.4
$ws.Range("M35").Value = -2847.4
$ws.Range("H45").Value = 3700
$ws.Range("I45").Value = 3700
$ws.Range("K45").Value = 3700
$ws.Range("M45").Value = -3323
$ws.Range("H102").Value = 4152.375
$ws.Range("I102").Value = 4152.375
$ws.Range("K102").Value = 4152.375
$ws.Range("M102").Value = -2530.375
$ws.Range("H122").Value = 1814.1428
$ws.Range("I122").Value = 1847.2106
$ws.Range("K122").Value = 5541.6318
$ws.Range("M122").Value = -3091.6318

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5666.6665
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377
$ws.Range("H89").Value = 5666.6665
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884
$ws.Range("H105").Value = 1682.7778
$ws.Range("I105").Value = 1743.125
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1743.125
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 3.875
$ws.Range("N105").Value = -4694

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8329.6
$ws.Range("I31").Value = 7613.7144
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 7613.7144
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -7318.7144
$ws.Range("N31").Value = -10590
$ws.Range("H34").Value = 8329.6
$ws.Range("I34").Value = 7613.7144
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 7613.7144
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -7411.7144
$ws.Range("N34").Value = -10404
$ws.Range("H44").Value = 11354.667
$ws.Range("I44").Value = 11354.667
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 11354.667
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -10912.667
$ws.Range("N44").ClearContents()
$ws.Range("H132").Value = 1677.9231
$ws.Range("I132").Value = 1518.5
$ws.Range("K132").Value = 4555.5
$ws.Range("M132").Value = -2025.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 14611.111
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 1500
$ws.Range("M35").Value = -1212
$ws.Range("H86").Value = 2379.5
$ws.Range("J86").Value = 3100
$ws.Range("L86").Value = 9300
$ws.Range("N86").Value = -11672
$ws.Range("H89").Value = 2379.5
$ws.Range("J89").Value = 3100
$ws.Range("L89").Value = 27900
$ws.Range("N89").Value = -39756
$ws.Range("H98").Value = 2182.7778
$ws.Range("I98").Value = 2481.8
$ws.Range("K98").Value = 7445.400000000001
$ws.Range("M98").Value = -5947.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H102").Value = 67782.836
$ws.Range("I102").Value = 80939.4
$ws.Range("K102").Value = 80939.4
$ws.Range("M102").Value = -79317.4
$ws.Range("H134").Value = 39665.332
$ws.Range("J134").Value = 39665.332
$ws.Range("L134").Value = 118995.996
$ws.Range("N134").Value = -124065.996

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1975.3334
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1975.3334
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1975.3334
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2565.3334
$ws.Range("H27").Value = 1975.3334
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1975.3334
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1975.3334
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2189.3334
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 28085.4
$ws.Range("I40").Value = 17606.125
$ws.Range("K40").Value = 17606.125
$ws.Range("M40").Value = -17470.125
$ws.Range("H47").Value = 3500
$ws.Range("J47").Value = 3500
$ws.Range("L47").Value = 3500
$ws.Range("N47").Value = -4480
$ws.Range("H52").Value = 3500
$ws.Range("J52").Value = 3500
$ws.Range("L52").Value = 3500
$ws.Range("N52").Value = -3966
$ws.Range("H58").Value = 26605
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H82").Value = 1472.5
$ws.Range("I82").Value = 1450
$ws.Range("K82").Value = 1450
$ws.Range("M82").Value = -1089
$ws.Range("H85").Value = 1472.5
$ws.Range("I85").Value = 1450
$ws.Range("K85").Value = 1450
$ws.Range("M85").Value = -202
$ws.Range("H132").Value = 999.4286
$ws.Range("I132").Value = 999.4286
$ws.Range("K132").Value = 2998.2858
$ws.Range("M132").Value = -468.2857999999997
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 11333.333
$ws.Range("I32").Value = 11333.333
$ws.Range("K32").Value = 11333.333
$ws.Range("M32").Value = -11016.333
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 11322.5
$ws.Range("J55").Value = 14993
$ws.Range("L55").Value = 14993
$ws.Range("N55").Value = -15547
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
